# Generate Report for Handoff
# Updates status for the f2a91ce9... file (row 3) from "Handed back: in sync with en-US"
# to "Ready for handoff" across the Overview, zh-cn, and de-de sheets, and refreshes the
# "Latest Handoff Datetime" timestamps for that file on the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("D2").Value = "2016-03-02 10:38:18"
$zhcn.Range("B3").Value = "Ready for handoff"
$zhcn.Range("D3").Value = "2016-03-02 10:38:18"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("D2").Value = "2016-03-02 10:38:29"
$dede.Range("B3").Value = "Ready for handoff"
$dede.Range("D3").Value = "2016-03-02 10:38:29"
